$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Technologies")
Write-Output ("col2=" + $ws.Columns.Item(2).ColumnWidth)
Write-Output ("col3=" + $ws.Columns.Item(3).ColumnWidth)
Write-Output ("col10=" + $ws.Columns.Item(10).ColumnWidth)
Write-Output ("col14=" + $ws.Columns.Item(14).ColumnWidth)
